$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'L1cam', 'Ptprz1', 'ECs', 3, 1, 6.603177, 19.809531, 0.5135477412645301, 0.5135477412645302, 2, 0.6666666666666666, 0.138766, 0.416298, 0.01356925767068476, 0.01356925767068476, 0.9162964595819999, 8.246668136238, 0.006968461627416556, 0.006968461627416558),
    @('ECs', 'L1cam', 'Ptprz1', 'FAPs', 3, 1, 6.603177, 19.809531, 0.5135477412645301, 0.5135477412645302, 2, 0.6666666666666666, 0.05416133333333333, 0.162484, 0.00529617548814441, 0.005296175488144411, 0.357636870556, 3.218731835004, 0.002719838959277132, 0.002719838959277133),
    @('ECs', 'L1cam', 'Ptprz1', 'MuSCs', 3, 1, 6.603177, 19.809531, 0.5135477412645301, 0.5135477412645302, 3, 1, 10.002366, 30.007098, 0.9780831152479456, 0.9780831152479456, 66.047393116782, 594.4265380510379, 0.5022923746045574, 0.5022923746045576),
    @('ECs', 'L1cam', 'Ptprz1', 'Resolving-Mac', 3, 1, 6.603177, 19.809531, 0.5135477412645301, 0.5135477412645302, 1, 0.3333333333333333, 0.03120566666666667, 0.09361700000000001, 0.003051451593225274, 0.003051451593225274, 0.206056540403, 1.854508863627, 0.001567066073278891, 0.001567066073278891),
    @('FAPs', 'L1cam', 'Ptprz1', 'ECs', 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.02944398858046029, 0.0294439885804603, 2, 0.6666666666666666, 0.138766, 0.416298, 0.01356925767068476, 0.01356925767068476, 0.05253537368466667, 0.472818363162, 0.0003995330679009652, 0.0003995330679009653),
    @('FAPs', 'L1cam', 'Ptprz1', 'FAPs', 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.02944398858046029, 0.0294439885804603, 2, 0.6666666666666666, 0.05416133333333333, 0.162484, 0.00529617548814441, 0.005296175488144411, 0.02050492113288889, 0.184544290196, 0.0001559405305930377, 0.0001559405305930378),
    @('FAPs', 'L1cam', 'Ptprz1', 'MuSCs', 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.02944398858046029, 0.0294439885804603, 3, 1, 10.002366, 30.007098, 0.9780831152479456, 0.9780831152479456, 3.786792409818, 34.081131688362, 0.02879866807610154, 0.02879866807610154),
    @('FAPs', 'L1cam', 'Ptprz1', 'Resolving-Mac', 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.02944398858046029, 0.0294439885804603, 1, 0.3333333333333333, 0.03120566666666667, 0.09361700000000001, 0.003051451593225274, 0.003051451593225274, 0.01181414294144444, 0.106327286473, 0.00008984690586475233, 0.00008984690586475233),
    @('MuSCs', 'L1cam', 'Ptprz1', 'ECs', 3, 1, 3.441487333333333, 10.324462, 0.2676541983690312, 0.2676541983690313, 2, 0.6666666666666666, 0.138766, 0.416298, 0.01356925767068476, 0.01356925767068476, 0.4775614312973334, 4.298052881676, 0.003631868784309956, 0.003631868784309958),
    @('MuSCs', 'L1cam', 'Ptprz1', 'FAPs', 3, 1, 3.441487333333333, 10.324462, 0.2676541983690312, 0.2676541983690313, 2, 0.6666666666666666, 0.05416133333333333, 0.162484, 0.00529617548814441, 0.005296175488144411, 0.1863955426231111, 1.677559883608, 0.001417543604701005, 0.001417543604701005),
    @('MuSCs', 'L1cam', 'Ptprz1', 'MuSCs', 3, 1, 3.441487333333333, 10.324462, 0.2676541983690312, 0.2676541983690313, 3, 1, 10.002366, 30.007098, 0.9780831152479456, 0.9780831152479456, 34.423015892364, 309.807143031276, 0.2617880521499736, 0.2617880521499737),
    @('MuSCs', 'L1cam', 'Ptprz1', 'Resolving-Mac', 3, 1, 3.441487333333333, 10.324462, 0.2676541983690312, 0.2676541983690313, 1, 0.3333333333333333, 0.03120566666666667, 0.09361700000000001, 0.003051451593225274, 0.003051451593225274, 0.1073939065615556, 0.9665451590540001, 0.0008167338300466137, 0.0008167338300466139),
    @('Resolving-Mac', 'L1cam', 'Ptprz1', 'ECs', 3, 1, 2.434707333333333, 7.304122, 0.1893540717859783, 0.1893540717859783, 2, 0.6666666666666666, 0.138766, 0.416298, 0.01356925767068476, 0.01356925767068476, 0.3378545978173333, 3.040691380356, 0.002569394191057278, 0.002569394191057279),
    @('Resolving-Mac', 'L1cam', 'Ptprz1', 'FAPs', 3, 1, 2.434707333333333, 7.304122, 0.1893540717859783, 0.1893540717859783, 2, 0.6666666666666666, 0.05416133333333333, 0.162484, 0.00529617548814441, 0.005296175488144411, 0.1318669954497778, 1.186802959048, 0.001002852393573235, 0.001002852393573236),
    @('Resolving-Mac', 'L1cam', 'Ptprz1', 'MuSCs', 3, 1, 2.434707333333333, 7.304122, 0.1893540717859783, 0.1893540717859783, 3, 1, 10.002366, 30.007098, 0.9780831152479456, 0.9780831152479456, 24.352833850884, 219.175504657956, 0.1852040204173128, 0.1852040204173128),
    @('Resolving-Mac', 'L1cam', 'Ptprz1', 'Resolving-Mac', 3, 1, 2.434707333333333, 7.304122, 0.1893540717859783, 0.1893540717859783, 1, 0.3333333333333333, 0.03120566666666667, 0.09361700000000001, 0.003051451593225274, 0.003051451593225274, 0.07597666547488889, 0.683789989274, 0.0005778047840350163, 0.0005778047840350164)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
